$d = $word.ActiveDocument

function Get-ListParagraphByText($doc, $exactText) {
    foreach ($para in $doc.Paragraphs) {
        $t = $para.Range.Text.TrimEnd([char]13, [char]7)
        if ($para.Style.NameLocal -eq "List Paragraph" -and $t -eq $exactText) {
            return $para
        }
    }
    return $null
}

# --- 1) "Either you rearrange the file structure ..." checklist item
#        (numId=2): add <w:strike/> to the paragraph mark rPr and to the
#        run rPr.
$text1 = "Either you rearrange the file structure to match the expected " + `
         "structure from the TEI and HTML files (but I assume you need to " + `
         "revise some XPaths as well) or you correct the XPaths in the XSL " + `
         "to reflect the existing file structure."
$p1 = Get-ListParagraphByText $d $text1
$p1.Range.Font.StrikeThrough = 1

# --- 2) "'cannot altogether map ... '" checklist item (numId=2): merge the
#        three runs ('  / text / ') into a single run and add <w:strike/>
#        to both the paragraph mark rPr and the run rPr.
$quote = [char]0x2018
$text2 = $quote + "cannot altogether map the contents of the TEI file with the contents of the html files. " + $quote
$p2 = Get-ListParagraphByText $d $text2
$r2 = $p2.Range
$r2.End = $r2.End - 1
$r2.Text = ""
$r2b = $p2.Range
$r2b.End = $r2b.End - 1
$r2b.InsertAfter($text2)
$p2.Range.Font.StrikeThrough = 1

# --- 3) "JPEGs needed? " checklist item: append a new run "- yes" (same
#        colour formatting as the rest of the paragraph, kept as its own
#        run).
$p3 = Get-ListParagraphByText $d "JPEGs needed? "
$r3 = $p3.Range
$r3.End = $r3.End - 1
$r3.Collapse(0)
$start3 = $r3.Start
$r3.InsertAfter("- yes")
$newRun = $d.Range($start3, $start3 + 5)
# toggling Bold on/off forces the inserted text into its own run while
# picking up the surrounding (black / theme text1) colour formatting
$newRun.Font.Bold = 1
$newRun.Font.Bold = 0

# --- 4) "Embed TIFFs with metadata" checklist item: add <w:strike/> to
#        both the paragraph mark rPr and the run rPr, keeping the existing
#        colour formatting.
$p4 = Get-ListParagraphByText $d "Embed TIFFs with metadata"
$p4.Range.Font.StrikeThrough = 1

Write-Output "done"
